$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 128
$ws.Range("I2").Value = 355
$ws.Range("J2").Value = 1531
$ws.Range("K2").Value = 4
$ws.Range("L2").Value = 413
$ws.Range("M2").Value = 22
$ws.Range("N2").Value = 264
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 7
$ws.Range("Q2").Value = 5
$ws.Range("R2").Value = 21
$ws.Range("S2").Value = 181
$ws.Range("T2").Value = 240
$ws.Range("U2").Value = 20
$ws.Range("V2").Value = 2350
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 2382
$ws.Range("Y2").Value = 4
$ws.Range("Z2").Value = 35
$ws.Range("AA2").Value = 15
